$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all the data values in the used range (A1:E6), while preserving
# cell styles/formatting - same effect as selecting the range and
# pressing Delete / Clear Contents.
$ws.Range("A1:E6").ClearContents()

# Move the active selection to E10, matching the saved cursor position.
$ws.Range("E10").Select()
